$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Year of Treatment" column (B). Everything to its right
# (C:Q) shifts left by one column to become B:P.
$ws.Columns("B:B").Delete()

# Append ".jamais.jamais" to every remaining header label in row 1,
# except the "Country" header in A1.
$headerRange = $ws.Range("B1:P1")
foreach ($cell in $headerRange.Cells) {
    $cell.Value = $cell.Value() + ".jamais.jamais"
}
